# Insert a new column at DE, shifting the existing "nom" (product name)
# column from DE->DF and "url_produit" from DF->DG for every row, then
# populate the newly-inserted DE column:
#   - DE1: new timestamp header "2026-02-01 16:15:30"
#   - DE2:DE80: numeric snapshot copied from the adjacent price column (DD)
#   - DE81:DE206: left blank (their DD counterpart is blank too)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything at/after column DE one column to the right.
$ws.Columns("DE:DE").Insert()

# New header cell for the inserted column (row 1).
$ws.Cells.Item(1, 109).Value = "2026-02-01 16:15:30"

# Populate the new column's numeric price snapshot for rows 2-80 (rows
# 81-206 keep the blank DD-style value that Insert() already left behind).
$ws.Cells.Item(2, 109).Value = 45.92
$ws.Cells.Item(3, 109).Value = 169.95
$ws.Cells.Item(4, 109).Value = 169.95
$ws.Cells.Item(5, 109).Value = 179.95
$ws.Cells.Item(6, 109).Value = 179.95
$ws.Cells.Item(7, 109).Value = 179.95
$ws.Cells.Item(8, 109).Value = 179.95
$ws.Cells.Item(9, 109).Value = 199.95
$ws.Cells.Item(10, 109).Value = 199.95
$ws.Cells.Item(11, 109).Value = 619
$ws.Cells.Item(12, 109).Value = 659
$ws.Cells.Item(13, 109).Value = 659
$ws.Cells.Item(14, 109).Value = 749
$ws.Cells.Item(15, 109).Value = 749
$ws.Cells.Item(16, 109).Value = 809
$ws.Cells.Item(17, 109).Value = 809
$ws.Cells.Item(18, 109).Value = 809
$ws.Cells.Item(19, 109).Value = 809
$ws.Cells.Item(20, 109).Value = 809
$ws.Cells.Item(21, 109).Value = 849
$ws.Cells.Item(22, 109).Value = 899
$ws.Cells.Item(23, 109).Value = 899
$ws.Cells.Item(24, 109).Value = 909
$ws.Cells.Item(25, 109).Value = 909
$ws.Cells.Item(26, 109).Value = 909
$ws.Cells.Item(27, 109).Value = 969
$ws.Cells.Item(28, 109).Value = 969
$ws.Cells.Item(29, 109).Value = 969
$ws.Cells.Item(30, 109).Value = 969
$ws.Cells.Item(31, 109).Value = 969
$ws.Cells.Item(32, 109).Value = 999
$ws.Cells.Item(33, 109).Value = 999
$ws.Cells.Item(34, 109).Value = 1039
$ws.Cells.Item(35, 109).Value = 1039
$ws.Cells.Item(36, 109).Value = 1079
$ws.Cells.Item(37, 109).Value = 1079
$ws.Cells.Item(38, 109).Value = 1079
$ws.Cells.Item(39, 109).Value = 1079
$ws.Cells.Item(40, 109).Value = 1099
$ws.Cells.Item(41, 109).Value = 1099
$ws.Cells.Item(42, 109).Value = 1199
$ws.Cells.Item(43, 109).Value = 1219
$ws.Cells.Item(44, 109).Value = 1219
$ws.Cells.Item(45, 109).Value = 1219
$ws.Cells.Item(46, 109).Value = 1219
$ws.Cells.Item(47, 109).Value = 1219
$ws.Cells.Item(48, 109).Value = 1229
$ws.Cells.Item(49, 109).Value = 1229
$ws.Cells.Item(50, 109).Value = 1249
$ws.Cells.Item(51, 109).Value = 1329
$ws.Cells.Item(52, 109).Value = 1329
$ws.Cells.Item(53, 109).Value = 1329
$ws.Cells.Item(54, 109).Value = 1329
$ws.Cells.Item(55, 109).Value = 1329
$ws.Cells.Item(56, 109).Value = 1329
$ws.Cells.Item(57, 109).Value = 1329
$ws.Cells.Item(58, 109).Value = 1349
$ws.Cells.Item(59, 109).Value = 1479
$ws.Cells.Item(60, 109).Value = 1479
$ws.Cells.Item(61, 109).Value = 1479
$ws.Cells.Item(62, 109).Value = 1579
$ws.Cells.Item(63, 109).Value = 1579
$ws.Cells.Item(64, 109).Value = 1579
$ws.Cells.Item(65, 109).Value = 1579
$ws.Cells.Item(66, 109).Value = 1579
$ws.Cells.Item(67, 109).Value = 1579
$ws.Cells.Item(68, 109).Value = 1579
$ws.Cells.Item(69, 109).Value = 1729
$ws.Cells.Item(70, 109).Value = 1729
$ws.Cells.Item(71, 109).Value = 1729
$ws.Cells.Item(72, 109).Value = 1829
$ws.Cells.Item(73, 109).Value = 1829
$ws.Cells.Item(74, 109).Value = 1829
$ws.Cells.Item(75, 109).Value = 1979
$ws.Cells.Item(76, 109).Value = 1979
$ws.Cells.Item(77, 109).Value = 1979
$ws.Cells.Item(78, 109).Value = 2479
$ws.Cells.Item(79, 109).Value = 2479
$ws.Cells.Item(80, 109).Value = 2479
